$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 127.61539
$ws.Range("I6").Value = 127
$ws.Range("J6").Value = 129.66667
$ws.Range("K6").Value = 381
$ws.Range("L6").Value = 389.00001
$ws.Range("M6").Value = -269
$ws.Range("N6").Value = -613.00001

$ws.Range("H40").Value = 2445.5
$ws.Range("I40").Value = 1591.1666
$ws.Range("K40").Value = 1591.1666
$ws.Range("M40").Value = -1416.1666

$ws.Range("H53").Value = 715.4286
$ws.Range("I53").Value = 676
$ws.Range("J53").Value = 754.8570999999999
$ws.Range("K53").Value = 676
$ws.Range("L53").Value = 754.8570999999999
$ws.Range("M53").Value = -39
$ws.Range("N53").Value = -2028.8571

$ws.Range("H62").Value = 4279.6
$ws.Range("I62").Value = 4275
$ws.Range("K62").Value = 4275
$ws.Range("M62").Value = -3651

$ws.Range("H65").Value = 4279.6
$ws.Range("I65").Value = 4275
$ws.Range("K65").Value = 21375
$ws.Range("M65").Value = -18255

$ws.Range("H74").Value = 5974.857
$ws.Range("I74").Value = 5137.3335
$ws.Range("K74").Value = 5137.3335
$ws.Range("M74").Value = -4201.3335

$ws.Range("H76").Value = 4933.3335
$ws.Range("J76").Value = 2400
$ws.Range("L76").Value = 2400
$ws.Range("N76").Value = -3030

$ws.Range("H77").Value = 5974.857
$ws.Range("I77").Value = 5137.3335
$ws.Range("K77").Value = 25686.6675
$ws.Range("M77").Value = -21006.6675

$ws.Range("H79").Value = 4933.3335
$ws.Range("J79").Value = 2400
$ws.Range("L79").Value = 2400
$ws.Range("N79").Value = -4584

$ws.Range("H86").Value = 10000
$ws.Range("I86").Value = 10000
$ws.Range("K86").Value = 10000
$ws.Range("M86").Value = -8877

$ws.Range("H89").Value = 10000
$ws.Range("I89").Value = 10000
$ws.Range("K89").Value = 50000
$ws.Range("M89").Value = -44384

$ws.Range("H92").Value = 516.2
$ws.Range("I92").Value = 516.2
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 516.2
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 731.8
$ws.Range("N92").ClearContents()

$ws.Range("H100").Value = 2271
$ws.Range("I100").Value = 2234.5
$ws.Range("K100").Value = 2234.5
$ws.Range("M100").Value = -1693.5

$ws.Range("H132").Value = 4302.1714
$ws.Range("I132").Value = 4512.7666
$ws.Range("J132").Value = 3038.6
$ws.Range("K132").Value = 13538.2998
$ws.Range("L132").Value = 9115.799999999999
$ws.Range("M132").Value = -11008.2998
$ws.Range("N132").Value = -14175.8

$ws.Range("H137").Value = 2571.3333
$ws.Range("I137").Value = 885
$ws.Range("J137").Value = 3998.2307
$ws.Range("K137").Value = 2655
$ws.Range("L137").Value = 11994.6921
$ws.Range("M137").Value = -105
$ws.Range("N137").Value = -17094.6921

$ws.Range("H138").Value = 2120.182
$ws.Range("J138").Value = 2374.6667
$ws.Range("L138").Value = 7124.000100000001
$ws.Range("N138").Value = -17404.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1485.5
$ws.Range("J88").Value = 1485.5
$ws.Range("L88").Value = 1485.5
$ws.Range("N88").Value = -2297.5

$ws.Range("H91").Value = 1485.5
$ws.Range("J91").Value = 1485.5
$ws.Range("L91").Value = 1485.5
$ws.Range("N91").Value = -4293.5

$ws.Range("H111").Value = 23248.25
$ws.Range("J111").Value = 23248.25
$ws.Range("L111").Value = 23248.25
$ws.Range("N111").Value = -31428.25

$ws.Range("H122").Value = 1763.75
$ws.Range("I122").Value = 1763.75
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5291.25
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2841.25
$ws.Range("N122").ClearContents()

$ws.Range("H133").Value = 66087
$ws.Range("J133").Value = 66087
$ws.Range("L133").Value = 66087
$ws.Range("N133").Value = -71147

$ws.Range("H134").Value = 75602.664
$ws.Range("J134").Value = 75602.664
$ws.Range("L134").Value = 75602.664
$ws.Range("N134").Value = -85742.664

$ws.Range("H140").Value = 99002.5
$ws.Range("J140").Value = 99002.5
$ws.Range("L140").Value = 99002.5
$ws.Range("N140").Value = -109362.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 939.44446
$ws.Range("I86").Value = 939.44446
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 939.44446
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 183.55554
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 939.44446
$ws.Range("I89").Value = 939.44446
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 4697.2223
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = 918.7776999999996
$ws.Range("N89").ClearContents()

$ws.Range("H94").Value = 1500
$ws.Range("I94").Value = 1500
$ws.Range("K94").Value = 1500
$ws.Range("M94").Value = -1049

$ws.Range("H105").Value = 5094.7
$ws.Range("I105").Value = 4996
$ws.Range("J105").Value = 5137
$ws.Range("K105").Value = 4996
$ws.Range("L105").Value = 5137
$ws.Range("M105").Value = -3249
$ws.Range("N105").Value = -8631

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2094
$ws.Range("I132").Value = 2094
$ws.Range("K132").Value = 6282
$ws.Range("M132").Value = -3752

$ws.Range("H135").Value = 118999
$ws.Range("J135").Value = 118999
$ws.Range("L135").Value = 118999
$ws.Range("N135").Value = -129139

$ws.Range("H140").Value = 123435
$ws.Range("J140").Value = 123435
$ws.Range("L140").Value = 123435
$ws.Range("N140").Value = -133795

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 9578
$ws.Range("I7").Value = 13871.182
$ws.Range("K7").Value = 41613.546
$ws.Range("M7").Value = -41501.546

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5182.5
$ws.Range("I70").Value = 1990.25
$ws.Range("K70").Value = 1990.25
$ws.Range("M70").Value = -1720.25

$ws.Range("H73").Value = 5182.5
$ws.Range("I73").Value = 1990.25
$ws.Range("K73").Value = 1990.25
$ws.Range("M73").Value = -1054.25

$ws.Range("H80").Value = 1142.8572
$ws.Range("I80").Value = 1000
$ws.Range("K80").Value = 1000
$ws.Range("M80").Value = -2

$ws.Range("H83").Value = 1142.8572
$ws.Range("I83").Value = 1000
$ws.Range("K83").Value = 5000
$ws.Range("M83").Value = -8

$ws.Range("H97").Value = 497.46155
$ws.Range("I97").Value = 463.1111
$ws.Range("K97").Value = 463.1111
$ws.Range("M97").Value = 32.88889999999998

$ws.Range("H102").Value = 2625
$ws.Range("I102").Value = 2750
$ws.Range("J102").Value = 2500
$ws.Range("K102").Value = 2750
$ws.Range("L102").Value = 2500
$ws.Range("M102").Value = -1128
$ws.Range("N102").Value = -5744

$ws.Range("H132").Value = 4257.3335
$ws.Range("I132").Value = 3608.8
$ws.Range("K132").Value = 10826.4
$ws.Range("M132").Value = -8296.400000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 6500
$ws.Range("I82").Value = 3500
$ws.Range("J82").Value = 11000
$ws.Range("K82").Value = 3500
$ws.Range("L82").Value = 11000
$ws.Range("M82").Value = -3139
$ws.Range("N82").Value = -11722

$ws.Range("H85").Value = 6500
$ws.Range("I85").Value = 3500
$ws.Range("J85").Value = 11000
$ws.Range("K85").Value = 3500
$ws.Range("L85").Value = 11000
$ws.Range("M85").Value = -2252
$ws.Range("N85").Value = -13496

$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

$ws.Range("H124").Value = 73328.664
$ws.Range("J124").Value = 73328.664
$ws.Range("L124").Value = 73328.664
$ws.Range("N124").Value = -83148.664

$ws.Range("H130").Value = 41649.332
$ws.Range("J130").Value = 41649.332
$ws.Range("L130").Value = 41649.332
$ws.Range("N130").Value = -51689.332

$ws.Range("H134").Value = 68999.28999999999
$ws.Range("J134").Value = 68999.28999999999
$ws.Range("L134").Value = 68999.28999999999
$ws.Range("N134").Value = -79139.28999999999

$ws.Range("H136").Value = 29720.285
$ws.Range("I136").Value = 10453.363
$ws.Range("J136").Value = 100365.664
$ws.Range("K136").Value = 31360.089
$ws.Range("L136").Value = 301096.992
$ws.Range("M136").Value = -28810.089
$ws.Range("N136").Value = -306196.992

$ws.Range("H137").Value = 70999
$ws.Range("I137").Value = 70999
$ws.Range("K137").Value = 70999
$ws.Range("M137").Value = -65899

$ws.Range("H139").Value = 99999.25
$ws.Range("J139").Value = 99999.25
$ws.Range("L139").Value = 99999.25
$ws.Range("N139").Value = -110279.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

$ws.Range("H132").Value = 1900.9584
$ws.Range("I132").Value = 1970.3478
$ws.Range("J132").Value = 305
$ws.Range("K132").Value = 5911.0434
$ws.Range("L132").Value = 915
$ws.Range("M132").Value = -3381.0434
$ws.Range("N132").Value = -5975

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws.Range("H136").Value = 3376.4285
$ws.Range("I136").Value = 3376.4285
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 10129.2855
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -7579.2855
$ws.Range("N136").ClearContents()

$ws.Range("H141").Value = 120000
$ws.Range("J141").Value = 120000
$ws.Range("L141").Value = 120000
